$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the weekly price-report data between row 3 (week of 2022-01-07)
# and row 4 (week of 2021-12-27): the rows were out of order and are
# being corrected back to chronological order.

# Row 3 <- (values that were in row 4)
$ws.Range("D3").Value = 44557
$ws.Range("J3").Value = 400
$ws.Range("K3").Value = 13000
$ws.Range("L3").Value = 14000
$ws.Range("M3").Value = 13500
$ws.Range("P3").Value = 750

# Row 4 <- (values that were in row 3)
$ws.Range("D4").Value = 44568
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 15000
$ws.Range("L4").Value = 16000
$ws.Range("M4").Value = 15500
$ws.Range("P4").Value = 861
